$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Day 8 (column M) was missing its completed-effort entries; bring M:O in
# line with the same row formatting used through column L (L6/L7), then
# fill in the day-8 numbers (M6 = Identify design patterns effort,
# M7 = Identify code smells effort). N/O stay blank, same as before.
$ws.Range("L6").Copy()
$ws.Range("M6:O6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("L7").Copy()
$ws.Range("M7:O7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("M6").Value = 4.0
$ws.Range("M7").Value = 0.0

# Remove the stray leftover note in L16 ("15  = 3 cada") that isn't part
# of the burndown table.
$ws.Range("L16").Clear()
